# Daily attendance processing - 2026-01-19 12:01:19
#
# The "Recorded By" column (G) lists the users/system accounts that
# recorded a session, as a comma-separated string. For a specific set of
# rows the last two names in that list were swapped (e.g.
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
#
# Rebuild each target cell by swapping the last two comma-separated
# tokens of its current value, so the edit is driven by the existing
# cell content rather than a hard-coded replacement string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,6,10,11,12,13,14,15,17,18,19,20,21,22,26,28,29,30,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,56,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = $cell.Value()
    $parts = @($current.Split(",") | ForEach-Object { $_.Trim() })

    $count = $parts.Count
    if ($count -ge 2) {
        $last = $parts[$count - 1]
        $secondLast = $parts[$count - 2]
        $parts[$count - 2] = $last
        $parts[$count - 1] = $secondLast
    }

    $cell.Value = [string]::Join(", ", $parts)
}
